# Apply updated cryptocurrency price/volume data to sheet1 of the workbook.
# All target cells hold text (inline string) values in the source file, including
# numeric-looking prices (e.g. "1.00", "0.0611") where trailing zeros / exact
# formatting must be preserved verbatim, so each cell is forced to Text format
# before assigning its string value (this stops Excel from "helpfully" coercing
# the text into a number and normalizing/rounding it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.789.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.671.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.129"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000201"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.153.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.626.13"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.681.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.33"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.70"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000117"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.69"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.60"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.97%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "534.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.30%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.50"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.12"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.84"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.11"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0611"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.641"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0257"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0989"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.82%  "
